{"js": "// Apply the commit's documented change:\n//   - Insert a new bullet \"Access denied, will show for several seconds\n//     after denied, once we refresh it will remain, until time has passed.\"\n//     right before the existing \"Test if we can log another workout...\" bullet.\n//   - Insert a new bullet \"Test location check, use google map to identify\n//     something within 10m of current location, change source code to\n//     reflect >20m, and see it get declined.\" right after the existing\n//     \"Click on total in gym to swap between, estimated total and actual\n//     total.\" bullet.\n// Both new bullets are plain top-level (\"ListParagraph\", numId 3, ilvl 0)\n// list items, matching their siblings in the same list.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfunction findByText(needle) {\n  for (const p of paragraphs.items) {\n    if ((p.text || \"\").indexOf(needle) !== -1) {\n      return p;\n    }\n  }\n  return null;\n}\n\n// 1) Insert the \"Access denied...\" bullet immediately before the\n//    \"Test if we can log another workout...\" paragraph.\nconst logWorkoutPara = findByText(\n  \"Test if we can log another workout after logging one today\"\n);\nif (!logWorkoutPara) {\n  throw new Error(\"Could not locate the 'Test if we can log another workout...' paragraph.\");\n}\nlogWorkoutPara.insertParagraph(\n  \"Access denied, will show for several seconds after denied, once we refresh it will remain, until time has passed.\",\n  \"Before\"\n);\n\n// 2) Insert the \"Test location check...\" bullet immediately after the\n//    \"Click on total in gym...\" paragraph, forcing it back to the\n//    top-level (ilvl 0) of the list (the anchor paragraph is a nested,\n//    ilvl 1, sub-bullet, so the inherited level must be corrected).\nconst clickTotalPara = findByText(\n  \"Click on total in gym to swap between, estimated total and actual total.\"\n);\nif (!clickTotalPara) {\n  throw new Error(\"Could not locate the 'Click on total in gym...' paragraph.\");\n}\nconst locationCheckPara = clickTotalPara.insertParagraph(\n  \"Test location check, use google map to identify something within 10m of current location, change source code to reflect >20m, and see it get declined.\",\n  \"After\"\n);\nawait context.sync();\n\nconst locationListItem = locationCheckPara.listItemOrNullObject;\nlocationListItem.load(\"isNullObject\");\nawait context.sync();\nif (!locationListItem.isNullObject) {\n  locationListItem.level = 0;\n}\n\nawait context.sync();\n", "ps1": "# Apply the commit's documented change:\n#   - Insert a new bullet \"Access denied, will show for several seconds\n#     after denied, once we refresh it will remain, until time has passed.\"\n#     right before the existing \"Test if we can log another workout...\" bullet.\n#   - Insert a new bullet \"Test location check, use google map to identify\n#     something within 10m of current location, change source code to\n#     reflect >20m, and see it get declined.\" right after the existing\n#     \"Click on total in gym to swap between, estimated total and actual\n#     total.\" bullet.\n# Both new bullets are plain top-level (\"ListParagraph\", numId 3, ilvl 0)\n# list items, matching their siblings in the same list.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndex($needle) {\n  for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*$needle*\") {\n      return $i\n    }\n  }\n  return -1\n}\n\n# 1) Insert the \"Access denied...\" bullet immediately before the\n#    \"Test if we can log another workout...\" paragraph.\n$logWorkoutIndex = Find-ParagraphIndex \"Test if we can log another workout after logging one today\"\nif ($logWorkoutIndex -eq -1) {\n  throw \"Could not locate the 'Test if we can log another workout...' paragraph.\"\n}\n\n$logWorkoutPara = $d.Paragraphs.Item($logWorkoutIndex)\n$logWorkoutPara.Range.InsertParagraphBefore()\n\n# The freshly inserted (empty) paragraph now sits at the same index the\n# target paragraph used to occupy; the target itself shifted down by one.\n$accessDeniedPara = $d.Paragraphs.Item($logWorkoutIndex)\n$accessDeniedPara.Range.Text = \"Access denied, will show for several seconds after denied, once we refresh it will remain, until time has passed.\"\n\n# 2) Insert the \"Test location check...\" bullet immediately after the\n#    \"Click on total in gym...\" paragraph, forcing it back to the\n#    top-level (ilvl 0 / ListLevelNumber 1) of the list (the anchor\n#    paragraph is a nested, ilvl 1, sub-bullet).\n$clickTotalIndex = Find-ParagraphIndex \"Click on total in gym to swap between, estimated total and actual total.\"\nif ($clickTotalIndex -eq -1) {\n  throw \"Could not locate the 'Click on total in gym...' paragraph.\"\n}\n\n$clickTotalPara = $d.Paragraphs.Item($clickTotalIndex)\n$clickTotalPara.Range.InsertParagraphAfter()\n\n$locationCheckPara = $d.Paragraphs.Item($clickTotalIndex + 1)\n$locationCheckPara.Range.Text = \"Test location check, use google map to identify something within 10m of current location, change source code to reflect >20m, and see it get declined.\"\n$locationCheckPara.Range.ListFormat.ListLevelNumber = 1\n"}
